# Daily refresh of the cryptos list (GitHub Actions scraper run).
#
# The sheet stores "Price" (column D) and "Volume(1h)" (column E) as plain
# text, not numbers -- prices keep a dotted thousands style (e.g.
# "63.095.59") and trailing zeros (e.g. "2.10") that a real Excel Number
# would normalize away, and the percentages keep their padded spacing
# (e.g. "  -2.43%  "). Column D entries that are NOT ambiguous (already
# contain more than one '.' or other non-numeric characters) are left as
# plain .Value assignments because Excel's COM layer already keeps those
# as text automatically. Column D entries that DO look like a plain
# number are round-tripped through a temporary "@" (Text) number format
# so Excel stores them as text instead of silently coercing them to a
# Double (which would drop the formatting), then the temporary format is
# cleared again so the cell keeps its original (default) style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "62.958.75"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").Value = "3.115.27"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.11%  "
Set-TextValue "D5" "593.51"
$ws.Range("E5").Value = "  -2.66%  "
Set-TextValue "D6" "135.88"
$ws.Range("E6").Value = "  -5.82%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "3.107.45"
$ws.Range("E8").Value = "  -1.13%  "
Set-TextValue "D9" "0.514"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("E10").Value = "  -4.29%  "
Set-TextValue "D11" "5.26"
$ws.Range("E11").Value = "  -2.73%  "
Set-TextValue "D12" "0.454"
$ws.Range("E12").Value = "  -2.95%  "
Set-TextValue "D13" "0.0000245"
$ws.Range("E13").Value = "  -5.29%  "
Set-TextValue "D14" "34.04"
$ws.Range("E14").Value = "  -4.08%  "
$ws.Range("D15").Value = "3.629.60"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("E16").Value = "  +2.55%  "
$ws.Range("D17").Value = "62.947.25"
$ws.Range("E17").Value = "  -2.21%  "
$ws.Range("D18").Value = "3.118.54"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("E19").Value = "  -3.14%  "
Set-TextValue "D20" "472.36"
$ws.Range("E20").Value = "  -1.84%  "
Set-TextValue "D21" "14.03"
$ws.Range("E21").Value = "  -4.79%  "
Set-TextValue "D22" "0.694"
$ws.Range("E22").Value = "  -3.23%  "
Set-TextValue "D23" "7.67"
$ws.Range("E23").Value = "  -0.75%  "
Set-TextValue "D24" "85.93"
$ws.Range("E24").Value = "  +0.47%  "
Set-TextValue "D25" "12.89"
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("E27").Value = "  -1.93%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D28" "7.89"
$ws.Range("E28").Value = "  -7.00%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D29" "6.92"
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("E31").Value = "  +0.08%  "
Set-TextValue "D32" "26.63"
$ws.Range("E32").Value = "  -1.43%  "
Set-TextValue "D33" "0.107"
$ws.Range("E33").Value = "  -6.77%  "
$ws.Range("E34").Value = "  -5.16%  "
Set-TextValue "D35" "1.07"
$ws.Range("E35").Value = "  -3.21%  "
Set-TextValue "D36" "5.78"
$ws.Range("E36").Value = "  -3.37%  "
Set-TextValue "D37" "52.03"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("D38").Value = "0.0₃0699"
$ws.Range("E38").Value = "  -9.17%  "
Set-TextValue "D39" "0.0387"
$ws.Range("E39").Value = "  -1.51%  "
Set-TextValue "D40" "418.29"
$ws.Range("E40").Value = "  -6.48%  "
Set-TextValue "D41" "8.17"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").Value = "2.897.91"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("E43").Value = "  -12.17%  "
$ws.Range("E44").Value = "  -6.04%  "
$ws.Range("E45").Value = "  +0.92%  "
Set-TextValue "D47" "2.10"
$ws.Range("E47").Value = "  -6.08%  "
Set-TextValue "D48" "25.42"
$ws.Range("E48").Value = "  -3.11%  "
Set-TextValue "D49" "0.112"
$ws.Range("E49").Value = "  -0.76%  "
Set-TextValue "D50" "2.24"
$ws.Range("E50").Value = "  -7.37%  "
Set-TextValue "D51" "119.44"
$ws.Range("E51").Value = "  -0.31%  "
